$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 33 (pushes old 33-37 down to 34-38)
[void]$ws.Rows("33").Insert()

# D31: stays empty, but gets a style (red, not bold, right aligned)
$ws.Range("D31").Font.Color = 255
$ws.Range("D31").HorizontalAlignment = -4152

# New row 33: note about seconds
$ws.Range("B33").Value = "All measurements are measured in seconds"
$ws.Range("A33:C33").Font.Bold = $true
$ws.Range("A33").Font.Color = 255
$ws.Range("A33").Font.Bold = $false
$ws.Range("A33").HorizontalAlignment = -4152

# A1: new label "Numbers are in seconds" - bold, italic, red font
$ws.Range("A1").Value = "Numbers are in seconds"
$ws.Range("A1").Font.Color = 255
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Italic = $true

[void]$ws.Range("B3").Select()

Write-Host "done"
